$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4100
$ws.Range("I29").Value = 3700
$ws.Range("K29").Value = 11100
$ws.Range("M29").Value = -10819
$ws.Range("H112").Value = 27212028
$ws.Range("J112").Value = 35715572
$ws.Range("L112").Value = 107146716
$ws.Range("N112").Value = -107148932
$ws.Range("H125").Value = 6683
$ws.Range("I125").Value = 11866
$ws.Range("K125").Value = 106794
$ws.Range("M125").Value = -104334
$ws.Range("H126").Value = 43280
$ws.Range("J126").Value = 43280
$ws.Range("L126").Value = 43280
$ws.Range("N126").Value = -53160
$ws.Range("H129").Value = 1144.8387
$ws.Range("I129").Value = 772.3
$ws.Range("K129").Value = 2316.9
$ws.Range("M129").Value = 2683.1
$ws.Range("H137").Value = 1469.1428
$ws.Range("I137").Value = 1376.3636
$ws.Range("K137").Value = 4129.0908
$ws.Range("M137").Value = -1579.0908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12861.848
$ws.Range("I32").Value = 8829.473
$ws.Range("J32").Value = 27378.4
$ws.Range("K32").Value = 8829.473
$ws.Range("L32").Value = 27378.4
$ws.Range("M32").Value = -8542.473
$ws.Range("N32").Value = -27952.4
$ws.Range("H61").Value = 306357.72
$ws.Range("I61").Value = 3176.7693
$ws.Range("J61").Value = 1432458.4
$ws.Range("K61").Value = 3176.7693
$ws.Range("L61").Value = 1432458.4
$ws.Range("M61").Value = -2964.7693
$ws.Range("N61").Value = -1432882.4
$ws.Range("H74").Value = 1459.1086
$ws.Range("I74").Value = 1190.65
$ws.Range("J74").Value = 1665.6154
$ws.Range("K74").Value = 1190.65
$ws.Range("L74").Value = 1665.6154
$ws.Range("M74").Value = -316.6500000000001
$ws.Range("N74").Value = -3413.6154
$ws.Range("H77").Value = 1459.1086
$ws.Range("I77").Value = 1190.65
$ws.Range("J77").Value = 1665.6154
$ws.Range("K77").Value = 5953.25
$ws.Range("L77").Value = 8328.076999999999
$ws.Range("M77").Value = -1585.25
$ws.Range("N77").Value = -17064.077
$ws.Range("H114").Value = 333353300
$ws.Range("J114").Value = 333353300
$ws.Range("L114").Value = 333353300
$ws.Range("N114").Value = -333361978
$ws.Range("H132").Value = 1616333.4
$ws.Range("I132").Value = 1776.4324
$ws.Range("J132").Value = 4005877.8
$ws.Range("K132").Value = 5329.2972
$ws.Range("L132").Value = 12017633.4
$ws.Range("M132").Value = -2799.2972
$ws.Range("N132").Value = -12022693.4
$ws.Range("H136").Value = 306357.72
$ws.Range("I136").Value = 3176.7693
$ws.Range("J136").Value = 1432458.4
$ws.Range("K136").Value = 9530.3079
$ws.Range("L136").Value = 4297375.199999999
$ws.Range("M136").Value = -6980.3079
$ws.Range("N136").Value = -4302475.199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1457.5518
$ws.Range("I94").Value = 1089.4762
$ws.Range("K94").Value = 1089.4762
$ws.Range("M94").Value = -638.4762000000001
$ws.Range("H105").Value = 3858.4285
$ws.Range("I105").Value = 1701.8
$ws.Range("J105").Value = 9250
$ws.Range("K105").Value = 1701.8
$ws.Range("L105").Value = 9250
$ws.Range("M105").Value = 45.20000000000005
$ws.Range("N105").Value = -12744
$ws.Range("H134").Value = 21212.223
$ws.Range("I134").Value = 2712.4893
$ws.Range("J134").Value = 145424.72
$ws.Range("K134").Value = 8137.467900000001
$ws.Range("L134").Value = 436274.16
$ws.Range("M134").Value = -5602.467900000001
$ws.Range("N134").Value = -441344.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6305.2363
$ws.Range("I31").Value = 1195.25
$ws.Range("J31").Value = 8401.641
$ws.Range("K31").Value = 1195.25
$ws.Range("L31").Value = 8401.641
$ws.Range("M31").Value = -900.25
$ws.Range("N31").Value = -8991.641
$ws.Range("H34").Value = 6305.2363
$ws.Range("I34").Value = 1195.25
$ws.Range("J34").Value = 8401.641
$ws.Range("K34").Value = 1195.25
$ws.Range("L34").Value = 8401.641
$ws.Range("M34").Value = -993.25
$ws.Range("N34").Value = -8805.641
$ws.Range("H43").Value = 41104.668
$ws.Range("J43").Value = 41104.668
$ws.Range("L43").Value = 41104.668
$ws.Range("N43").Value = -41472.668
$ws.Range("H58").Value = 835429.0600000001
$ws.Range("I58").Value = 2286.182
$ws.Range("K58").Value = 2286.182
$ws.Range("M58").Value = -2083.182
$ws.Range("H99").Value = 5954077
$ws.Range("I99").Value = 1433.7222
$ws.Range("K99").Value = 1433.7222
$ws.Range("M99").Value = 64.27780000000007
$ws.Range("H101").Value = 41104.668
$ws.Range("J101").Value = 41104.668
$ws.Range("L101").Value = 41104.668
$ws.Range("N101").Value = -47594.668
$ws.Range("H126").Value = 5954077
$ws.Range("I126").Value = 1433.7222
$ws.Range("K126").Value = 4301.1666
$ws.Range("M126").Value = -1831.1666
$ws.Range("H132").Value = 1627.5834
$ws.Range("I132").Value = 954.25
$ws.Range("K132").Value = 2862.75
$ws.Range("M132").Value = -332.75
$ws.Range("H136").Value = 835429.0600000001
$ws.Range("I136").Value = 2286.182
$ws.Range("K136").Value = 6858.545999999999
$ws.Range("M136").Value = -4308.545999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 448.61905
$ws.Range("I107").Value = 451.18182
$ws.Range("J107").Value = 445.8
$ws.Range("K107").Value = 1353.54546
$ws.Range("L107").Value = 1337.4
$ws.Range("M107").Value = 566.45454
$ws.Range("N107").Value = -5177.4
$ws.Range("H131").Value = 1613985.5
$ws.Range("I131").Value = 4762416
$ws.Range("J131").Value = 1374.7805
$ws.Range("K131").Value = 14287248
$ws.Range("L131").Value = 4124.3415
$ws.Range("M131").Value = -14282208
$ws.Range("N131").Value = -14204.3415
$ws.Range("H140").Value = 3601.4092
$ws.Range("I140").Value = 3601.4092
$ws.Range("K140").Value = 10804.2276
$ws.Range("M140").Value = -5624.2276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8587.75
$ws.Range("J92").Value = 8587.75
$ws.Range("L92").Value = 8587.75
$ws.Range("N92").Value = -12331.75
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H113").Value = 142858670
$ws.Range("I113").Value = 1000000000
$ws.Range("K113").Value = 1000000000
$ws.Range("M113").Value = -999997830
$ws.Range("H122").Value = 43829070
$ws.Range("I122").Value = 53242360
$ws.Range("J122").Value = 25002490
$ws.Range("K122").Value = 159727080
$ws.Range("L122").Value = 75007470
$ws.Range("M122").Value = -159724630
$ws.Range("N122").Value = -75012370
$ws.Range("H126").Value = 13930.3125
$ws.Range("I126").Value = 17407.5
$ws.Range("J126").Value = 3498.75
$ws.Range("K126").Value = 52222.5
$ws.Range("L126").Value = 10496.25
$ws.Range("M126").Value = -49752.5
$ws.Range("N126").Value = -15436.25
$ws.Range("H132").Value = 12434
$ws.Range("I132").Value = 9725.916999999999
$ws.Range("J132").Value = 23266.334
$ws.Range("K132").Value = 29177.751
$ws.Range("L132").Value = 69799.00199999999
$ws.Range("M132").Value = -26647.751
$ws.Range("N132").Value = -74859.00199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 900.3333
$ws.Range("I16").Value = 900.3333
$ws.Range("K16").Value = 900.3333
$ws.Range("M16").Value = -730.3333
$ws.Range("H40").Value = 111113656
$ws.Range("I40").Value = 142859550
$ws.Range("K40").Value = 142859550
$ws.Range("M40").Value = -142859414
$ws.Range("H82").Value = 111598.5
$ws.Range("I82").Value = 1776.8
$ws.Range("J82").Value = 221420.2
$ws.Range("K82").Value = 1776.8
$ws.Range("L82").Value = 221420.2
$ws.Range("M82").Value = -1415.8
$ws.Range("N82").Value = -222142.2
$ws.Range("H85").Value = 111598.5
$ws.Range("I85").Value = 1776.8
$ws.Range("J85").Value = 221420.2
$ws.Range("K85").Value = 1776.8
$ws.Range("L85").Value = 221420.2
$ws.Range("M85").Value = -528.8
$ws.Range("N85").Value = -223916.2
$ws.Range("H103").Value = 42534
$ws.Range("J103").Value = 42534
$ws.Range("L103").Value = 42534
$ws.Range("N103").Value = -44878

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 694.4516
$ws.Range("I126").Value = 632
$ws.Range("J126").Value = 908.5714
$ws.Range("K126").Value = 1896
$ws.Range("L126").Value = 2725.7142
$ws.Range("M126").Value = 574
$ws.Range("N126").Value = -7665.7142
$ws.Range("H132").Value = 1991.4348
$ws.Range("I132").Value = 793.3333
$ws.Range("J132").Value = 4237.875
$ws.Range("K132").Value = 2379.9999
$ws.Range("L132").Value = 12713.625
$ws.Range("M132").Value = 150.0001000000002
$ws.Range("N132").Value = -17773.625
$ws.Range("H136").Value = 5420
$ws.Range("I136").Value = 6594.4443
$ws.Range("J136").Value = 3658.3333
$ws.Range("K136").Value = 19783.3329
$ws.Range("L136").Value = 10974.9999
$ws.Range("M136").Value = -17233.3329
$ws.Range("N136").Value = -16074.9999
